$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.965.37'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.305.43'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '185.65'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '575.87'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -0.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.66'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.412'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.875.15'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.39'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.250.41'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.94%  '

$ws.Range("E16").Value = '  -0.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.281.49'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '442.82'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +9.90%  '

$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.50'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.60%  '

$ws.Range("E21").Value = '  +2.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.07'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.14%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.512'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.439.34'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.91%  '

$ws.Range("E26").Value = '  +1.40%  '

$ws.Range("E27").Value = '  -0.37%  '

$ws.Range("E28").Value = '  -4.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.55%  '

$ws.Range("E30").Value = '  +1.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.81'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.33'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.56%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("E34").Value = '  -1.10%  '

$ws.Range("E35").Value = '  -1.87%  '

$ws.Range("E36").Value = '  +4.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.09'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.49'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.84'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("E40").Value = '  -1.57%  '

$ws.Range("E41").Value = '  -0.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.736.09'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.29'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.26'
$ws.Range("D44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0671'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.73'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.54%  '

$ws.Range("E47").Value = '  -1.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '328.81'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.990'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.25'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.85%  '
